$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '35.575.22'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +1.69%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.907.10'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +3.48%  '
$ws.Range('E4').Value = '  +0.66%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '245.92'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +5.54%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.632'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +2.10%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '42.40'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +2.31%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.338'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +3.34%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.0705'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.75%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0998'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.74%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '2.184.29'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +3.55%  '
$ws.Range('E13').Value = '  +8.73%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '1.906.29'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +3.57%  '
$ws.Range('E15').Value = '  +2.92%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '4.85'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +3.39%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '35.578.53'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.66%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '71.96'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.84%  '
$ws.Range('E19').Value = '  +2.56%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '243.73'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.50%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '12.44'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +2.26%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.92'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +3.41%  '
$ws.Range('E23').Value = '  +0.67%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.28'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -1.04%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '172.03'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.50%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.20'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +26.69%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.54'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +8.53%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '17.95'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('E29').Value = '  +1.11%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.979'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +29.61%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.10'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +3.72%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.0564'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +2.06%  '
$ws.Range('E33').Value = '  +0.87%  '
$ws.Range('E34').Value = '  +5.61%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.74'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +6.03%  '
$ws.Range('E36').Value = '  +3.24%  '
$ws.Range('E37').Value = '  +4.57%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '54.12'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +58.12%  '
$ws.Range('E39').Value = '  +4.82%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.0205'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +3.21%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '91.43'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.72%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.355.87'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +0.94%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '15.41'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +6.15%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.0592'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +11.91%  '
$ws.Range('E45').Value = '  +3.69%  '
$ws.Range('B46').Value = 'Gas'
$ws.Range('C46').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '12.71'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +7.44%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.44'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.48%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  +5.64%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.093.71'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +3.54%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.0688'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +2.45%  '
